$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "'61.743.97"
$ws.Range("E2").Value2 = "  -7.73%  "

# Row 3
$ws.Range("D3").Value2 = "'2.908.52"
$ws.Range("E3").Value2 = "  -9.74%  "

# Row 4
$ws.Range("D4").Value2 = "'1.00"
$ws.Range("E4").Value2 = "  +0.02%  "

# Row 5
$ws.Range("D5").Value2 = "'524.20"
$ws.Range("E5").Value2 = "  -11.25%  "

# Row 6
$ws.Range("D6").Value2 = "'123.92"
$ws.Range("E6").Value2 = "  -18.36%  "

# Row 7
$ws.Range("E7").Value2 = "  +0.06%  "

# Row 8
$ws.Range("D8").Value2 = "'2.878.20"
$ws.Range("E8").Value2 = "  -10.39%  "

# Row 9
$ws.Range("D9").Value2 = "'0.433"
$ws.Range("E9").Value2 = "  -20.12%  "

# Row 10
$ws.Range("B10").Value2 = "Toncoin"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value2 = "'5.77"
$ws.Range("E10").Value2 = "  -10.41%  "

# Row 11
$ws.Range("B11").Value2 = "Dogecoin"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value2 = "'0.138"
$ws.Range("E11").Value2 = "  -19.41%  "

# Row 12
$ws.Range("D12").Value2 = "'0.413"
$ws.Range("E12").Value2 = "  -16.18%  "

# Row 13
$ws.Range("D13").Value2 = "'0.0000195"
$ws.Range("E13").Value2 = "  -19.65%  "

# Row 14
$ws.Range("D14").Value2 = "'30.40"
$ws.Range("E14").Value2 = "  -21.84%  "

# Row 15
$ws.Range("D15").Value2 = "'3.413.57"
$ws.Range("E15").Value2 = "  -8.97%  "

# Row 16
$ws.Range("D16").Value2 = "'61.700.47"
$ws.Range("E16").Value2 = "  -7.93%  "

# Row 17
$ws.Range("E17").Value2 = "  -4.68%  "

# Row 18
$ws.Range("D18").Value2 = "'2.929.38"
$ws.Range("E18").Value2 = "  -9.43%  "

# Row 19
$ws.Range("D19").Value2 = "'460.93"
$ws.Range("E19").Value2 = "  -12.94%  "

# Row 20
$ws.Range("D20").Value2 = "'5.92"
$ws.Range("E20").Value2 = "  -16.60%  "

# Row 21
$ws.Range("D21").Value2 = "'12.19"
$ws.Range("E21").Value2 = "  -17.81%  "

# Row 22
$ws.Range("D22").Value2 = "'0.608"
$ws.Range("E22").Value2 = "  -19.60%  "

# Row 23
$ws.Range("D23").Value2 = "'6.22"
$ws.Range("E23").Value2 = "  -21.42%  "

# Row 24
$ws.Range("D24").Value2 = "'72.67"
$ws.Range("E24").Value2 = "  -14.98%  "

# Row 25
$ws.Range("D25").Value2 = "'0.997"
$ws.Range("E25").Value2 = "  -0.12%  "

# Row 26
$ws.Range("D26").Value2 = "'11.29"
$ws.Range("E26").Value2 = "  -18.11%  "

# Row 27
$ws.Range("D27").Value2 = "'2.57"
$ws.Range("E27").Value2 = "  -19.22%  "

# Row 28
$ws.Range("D28").Value2 = "'6.62"
$ws.Range("E28").Value2 = "  -17.64%  "

# Row 29
$ws.Range("D29").Value2 = "'1.75"
$ws.Range("E29").Value2 = "  -19.41%  "

# Row 30
$ws.Range("D30").Value2 = "'23.74"
$ws.Range("E30").Value2 = "  -18.35%  "

# Row 31
$ws.Range("D31").Value2 = "'1.04"
$ws.Range("E31").Value2 = "  -9.70%  "

# Row 32
$ws.Range("D32").Value2 = "'1.00"
$ws.Range("E32").Value2 = "  -0.51%  "

# Row 33
$ws.Range("D33").Value2 = "'2.19"
$ws.Range("E33").Value2 = "  -17.90%  "

# Row 34
$ws.Range("B34").Value2 = "Bittensor"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value2 = "'460.65"
$ws.Range("E34").Value2 = "  -15.63%  "

# Row 35
$ws.Range("B35").Value2 = "OKB"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value2 = "'50.61"
$ws.Range("E35").Value2 = "  -5.52%  "

# Row 36
$ws.Range("D36").Value2 = "'5.22"
$ws.Range("E36").Value2 = "  -18.92%  "

# Row 37
$ws.Range("D37").Value2 = "'4.51"
$ws.Range("E37").Value2 = "  -21.32%  "

# Row 38
$ws.Range("D38").Value2 = "'0.0372"
$ws.Range("E38").Value2 = "  -12.51%  "

# Row 39
$ws.Range("D39").Value2 = "'0.0730"
$ws.Range("E39").Value2 = "  -15.26%  "

# Row 40
$ws.Range("D40").Value2 = "'0.107"
$ws.Range("E40").Value2 = "  -14.25%  "

# Row 41
$ws.Range("D41").Value2 = "'7.45"
$ws.Range("E41").Value2 = "  -19.57%  "

# Row 42
$ws.Range("D42").Value2 = "'2.571.66"
$ws.Range("E42").Value2 = "  -11.90%  "

# Row 43
$ws.Range("E43").Value2 = "  -0.25%  "

# Row 44
$ws.Range("D44").Value2 = "'2.13"
$ws.Range("E44").Value2 = "  -19.72%  "

# Row 45
$ws.Range("B45").Value2 = "TheGraph"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value2 = "'0.214"
$ws.Range("E45").Value2 = "  -18.75%  "

# Row 46
$ws.Range("B46").Value2 = "Monero"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value2 = "'109.78"
$ws.Range("E46").Value2 = "  -6.97%  "

# Row 47
$ws.Range("D47").Value2 = "'0.0967"
$ws.Range("E47").Value2 = "  -15.68%  "

# Row 48
$ws.Range("D48").Value2 = "'1.72"
$ws.Range("E48").Value2 = "  -19.86%  "

# Row 49
$ws.Range("B49").Value2 = "PEPE"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value2 = "'0.0₃0445"
$ws.Range("E49").Value2 = "  -23.50%  "

# Row 50
$ws.Range("B50").Value2 = "BitgetToken"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D50").Value2 = "'1.15"
$ws.Range("E50").Value2 = "  -8.08%  "

# Row 51
$ws.Range("B51").Value2 = "InjectiveProtocol"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value2 = "'20.74"
$ws.Range("E51").Value2 = "  -21.45%  "

